{"js": "// Prepend \"Ng\u00e0y chuy\u1ec3n \u0111\u1ed5i: \" (as separate word/space runs with spell-check\n// markers, matching Word's own AutoCorrect/proofing run-splitting) in front\n// of the \"<conversionDate>\" placeholder inside the signature table, leaving\n// the placeholder run's own formatting (w:i) untouched.\nconst results = context.document.body.search(\"<conversionDate>\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<conversionDate>' placeholder in the document body.\");\n}\n\nconst target = results.items[0];\n\n// Build the replacement as raw OOXML so the inserted runs carry the exact\n// formatting/proofing markup used by Word (iCs + lang, split on word\n// boundaries with proofErr spellStart/spellEnd wrappers), while the\n// existing \"<conversionDate>\" run (kept at the end, still italic) is\n// preserved as part of the same paragraph.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Ng\u00e0y</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>chuy\u1ec3n</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>\u0111\u1ed5i</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>\n<w:r><w:rPr><w:i/></w:rPr><w:t>&lt;conversionDate&gt;</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Prepend \"Ng\u00e0y chuy\u1ec3n \u0111\u1ed5i: \" (split into the same word/space runs with\n# proofing spell-check markers that Word itself would produce) in front of\n# the \"<conversionDate>\" placeholder inside the signature table, while\n# leaving the placeholder run's own formatting (w:i) untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"<conversionDate>\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find '<conversionDate>' placeholder in the document.\"\n}\n\n# Rebuild the whole paragraph (same paraId/rsid attributes and pPr as the\n# original) with the new runs inserted ahead of the existing placeholder\n# run, by replacing the found range's containing markup via raw OOXML.\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n<w:body>\n<w:p w14:paraId=\"783E2BAA\" w14:textId=\"32F5FAFB\" w:rsidR=\"00B5588D\" w:rsidRPr=\"007B7260\" w:rsidRDefault=\"00101574\" w:rsidP=\"007B7260\">\n<w:pPr><w:pStyle w:val=\"TableParagraph\"/><w:jc w:val=\"center\"/><w:rPr><w:i/></w:rPr></w:pPr>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Ng\u00e0y</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>chuy\u1ec3n</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t>\u0111\u1ed5i</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:rPr><w:iCs/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>\n<w:r><w:rPr><w:i/></w:rPr><w:t>&lt;conversionDate&gt;</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$range.InsertXML($ooxml)\n"}
